$d = $word.ActiveDocument

# The "Hacer que el tamaño de la ventana no sea modificable" to-do item
# (window resizing) was folded into the new separate component class, so
# the task line - and the now-unused trailing blank paragraph after it -
# are removed from the list.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Hacer que el tamaño de la ventana no sea modificable*") {
        $target = $p
        break
    }
}

$next = $target.Next()
$rng = $d.Range($target.Range.Start, $next.Range.End)
$rng.Delete()
